$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the closing-signature block from rows 23/24 down to rows 31/32
#    (it gets pushed down by the newly inserted worker rows). Do this BEFORE
#    the data table rows below are touched, since it reads their format.
# ---------------------------------------------------------------------------
$ws.Range("B23:C24").Copy()
$ws.Range("B31:C32").PasteSpecial(-4122)
$ws.Range("H23:J24").Copy()
$ws.Range("H31:J32").PasteSpecial(-4122)

$ws.Range("B23:C23").UnMerge()
$ws.Range("B24:C24").UnMerge()
$ws.Range("H23:J23").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B31:C31").Merge()
$ws.Range("B32:C32").Merge()
$ws.Range("H31:J31").Merge()
$ws.Range("H32:J32").Merge()

$ws.Range("B31").Value = "___________________________________"
$ws.Range("H31").Value = "___________________________________"
$ws.Range("B32").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H32").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# ---------------------------------------------------------------------------
# 2. Prepare the worker-rows block (16-26). Row 26 keeps the "closing" look
#    (thicker outer border) that row 18 used to have; rows 16-25 use the
#    "standard" look that rows 16/17 already have.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J25").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Fill in the new worker data (Tipo Doc, N Doc, Nombre, Periodo, Valor
#    Mora, Salario Basico). Columns B-E are textual (document numbers /
#    periods must stay text, not get auto-converted to numbers), so force
#    Text format before assigning them.
# ---------------------------------------------------------------------------
$rows = @(
  @{r=16; tipo="CC"; doc="9298356";     nombre="WALBERTO ROMERO GONZALEZ";          periodo="2507"; mora=56940;  salario=877803},
  @{r=17; tipo="CC"; doc="22801874";    nombre="LIDA MARGARITA GARCIA ARIAS";       periodo="2507"; mora=208926; salario=5223150},
  @{r=18; tipo="CC"; doc="1001972388";  nombre="LESLY RODRIGUEZ SALCEDO";           periodo="2507"; mora=72270;  salario=737717},
  @{r=19; tipo="CC"; doc="1050958276";  nombre="MIGUEL ANTONIO ROMERO TEHERAN";     periodo="2507"; mora=56940;  salario=908526},
  @{r=20; tipo="CC"; doc="79774845";    nombre="SANDRO RODRIGUEZ ZULUAGA";          periodo="2506"; mora=9490;   salario=1423500},
  @{r=21; tipo="CC"; doc="1143348955";  nombre="YESID ALEXANDER MONTERREY PORTO";   periodo="2507"; mora=56940;  salario=1423500},
  @{r=22; tipo="CC"; doc="1143348955";  nombre="YESID ALEXANDER MONTERREY PORTO";   periodo="2506"; mora=9490;   salario=1423500},
  @{r=23; tipo="CC"; doc="1047454082";  nombre="NAYBER ENRIQUE TAPIA PIMIENTA";     periodo="2507"; mora=56940;  salario=1423500},
  @{r=24; tipo="CC"; doc="1047442373";  nombre="SUAD HELENA ROA ANGULO";            periodo="2507"; mora=56940;  salario=1423500},
  @{r=25; tipo="CC"; doc="1002444328";  nombre="CRISTIAN DAVID ROMERO MEDINA";      periodo="2507"; mora=56940;  salario=1423500},
  @{r=26; tipo="CC"; doc="1019060050";  nombre="EDWAR ALEXANDER MUÑOZ REYES";       periodo="2506"; mora=9490;   salario=1423500}
)

foreach ($row in $rows) {
  $r = $row.r

  $ws.Cells.Item($r, 2).Value = $row.tipo

  $ws.Cells.Item($r, 3).NumberFormat = "@"
  $ws.Cells.Item($r, 3).Value = $row.doc

  $ws.Cells.Item($r, 4).Value = $row.nombre

  $ws.Cells.Item($r, 5).NumberFormat = "@"
  $ws.Cells.Item($r, 5).Value = $row.periodo

  $ws.Cells.Item($r, 6).Value = $row.mora
  $ws.Cells.Item($r, 7).Value = $row.salario
}

# ---------------------------------------------------------------------------
# 4. Update the summary header fields: total overdue amount, worker count,
#    period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 651306
$ws.Range("C13").Value = 10
$ws.Range("F13").Value = 2

# ---------------------------------------------------------------------------
# 5. Column D (Nombre Trabajador) widens to fit the longest new name.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 36.90625
